$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (D:E) shifting existing D:K to F:M, only within the used row range
$ws.Range("D5:E102").Insert(-4161)

# Copy formatting (number formats/styles) from the old column positions (now F:G) into the new D:E columns
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 1132700
$ws.Range("E8").Value2 = 775100
$ws.Range("D9").Value2 = 897900
$ws.Range("E9").Value2 = 609900
$ws.Range("D10").Value2 = 234800
$ws.Range("E10").Value2 = 165200
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("D17").Value2 = 1000100
$ws.Range("E17").Value2 = 693000
$ws.Range("D18").Value2 = 132600
$ws.Range("E18").Value2 = 82100
$ws.Range("D20").Value2 = 3500
$ws.Range("E20").Value2 = 1500
$ws.Range("D21").Value2 = 145600
$ws.Range("E21").Value2 = 90600
$ws.Range("D22").Value2 = 0
$ws.Range("E22").Value2 = 0
$ws.Range("D23").Value2 = 136100
$ws.Range("E23").Value2 = 83600
$ws.Range("D24").Value2 = 35800
$ws.Range("E24").Value2 = 19700
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 100200
$ws.Range("E26").Value2 = 64000
$ws.Range("D27").Value2 = 98600
$ws.Range("E27").Value2 = 64000
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 700
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = -3500
$ws.Range("E32").Value2 = -1500
$ws.Range("D33").Value2 = 99400
$ws.Range("E33").Value2 = 64000
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 99400
$ws.Range("E35").Value2 = 64000
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 277700
$ws.Range("E41").Value2 = 83100
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("D43").Value2 = 51600
$ws.Range("E43").Value2 = 85000
$ws.Range("D44").Value2 = 3216100
$ws.Range("E44").Value2 = 3377700
$ws.Range("D45").Value2 = 32000
$ws.Range("E45").Value2 = 27000
$ws.Range("D46").Value2 = 0
$ws.Range("E46").Value2 = 0
$ws.Range("D47").Value2 = 5400
$ws.Range("E47").Value2 = 4300
$ws.Range("D48").Value2 = 0
$ws.Range("E48").Value2 = 0
$ws.Range("D49").Value2 = 160400
$ws.Range("E49").Value2 = 160600
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 70200
$ws.Range("E52").Value2 = 61800
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 3884200
$ws.Range("E54").Value2 = 3877100
$ws.Range("D57").Value2 = 81300
$ws.Range("E57").Value2 = 83700
$ws.Range("D58").Value2 = 0
$ws.Range("E58").Value2 = 0
$ws.Range("D59").Value2 = 297600
$ws.Range("E59").Value2 = 293400
$ws.Range("D60").Value2 = 0
$ws.Range("E60").Value2 = 0
$ws.Range("D61").Value2 = 1410800
$ws.Range("E61").Value2 = 1519200
$ws.Range("D62").Value2 = 0
$ws.Range("E62").Value2 = 0
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 1827300
$ws.Range("E66").Value2 = 1916700
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 1396800
$ws.Range("E72").Value2 = 1297400
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 2056900
$ws.Range("E76").Value2 = 1960400
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 99400
$ws.Range("E81").Value2 = 64000
$ws.Range("D83").Value2 = 9500
$ws.Range("E83").Value2 = 7000
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 382800
$ws.Range("E89").Value2 = -72000
$ws.Range("D91").Value2 = -7100
$ws.Range("E91").Value2 = -8900
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -69100
$ws.Range("E94").Value2 = -9500
$ws.Range("D96").Value2 = 0
$ws.Range("E96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -119100
$ws.Range("E100").Value2 = -75300
$ws.Range("D101").Value2 = 0
$ws.Range("E101").Value2 = 0
$ws.Range("D102").Value2 = 194600
$ws.Range("E102").Value2 = -156800

Write-Host "edit complete"
